$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042453823061956
$ws.Range("D2").Value = 1.050059299706727
$ws.Range("E2").Value = 1.056309181877289
$ws.Range("F2").Value = 1.062944626805408
$ws.Range("I2").Value = 1.03940610317222
$ws.Range("J2").Value = 1.047529560930557
$ws.Range("K2").Value = 1.052814629854615
$ws.Range("L2").Value = 1.059047246760224
$ws.Range("M2").Value = 1.065664601864454
$ws.Range("N2").Value = 1.0196536817771
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.043635600619232
$ws.Range("D3").Value = 1.050974513174485
$ws.Range("E3").Value = 1.057342277321322
$ws.Range("F3").Value = 1.063990447109192
$ws.Range("I3").Value = 1.0396440004701
$ws.Range("J3").Value = 1.048356628932064
$ws.Range("K3").Value = 1.053541948625701
$ws.Range("L3").Value = 1.059893396887044
$ws.Range("M3").Value = 1.066524756850403
$ws.Range("N3").Value = 1.019935004194602
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044400224528282
$ws.Range("D4").Value = 1.051566530633798
$ws.Range("E4").Value = 1.0580114406879
$ws.Range("F4").Value = 1.064667483877108
$ws.Range("I4").Value = 1.039796534382274
$ws.Range("J4").Value = 1.048891207229141
$ws.Range("K4").Value = 1.054011756249492
$ws.Range("L4").Value = 1.060440986180277
$ws.Range("M4").Value = 1.067081049616671
$ws.Range("N4").Value = 1.020116671886963
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044721658085566
$ws.Range("D5").Value = 1.051815370156864
$ws.Range("E5").Value = 1.05829292009656
$ws.Range("F5").Value = 1.064952187132097
$ws.Range("I5").Value = 1.039860324112157
$ws.Range("J5").Value = 1.04911580319181
$ws.Range("K5").Value = 1.054209068188491
$ws.Range("L5").Value = 1.060671210367289
$ws.Range("M5").Value = 1.067314847012016
$ws.Range("N5").Value = 1.020192957249295
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.044775627390731
$ws.Range("D6").Value = 1.05185714879474
$ws.Range("E6").Value = 1.058340191292936
$ws.Range("F6").Value = 1.064999994580218
$ws.Range("I6").Value = 1.039871015013235
$ws.Range("J6").Value = 1.049153505602756
$ws.Range("K6").Value = 1.054242186292058
$ws.Range("L6").Value = 1.060709867057177
$ws.Range("M6").Value = 1.06735409863543
$ws.Range("N6").Value = 1.020205760754145
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04440451959393
$ws.Range("D7").Value = 1.051569855814859
$ws.Range("E7").Value = 1.058015201188588
$ws.Range("F7").Value = 1.064671287793231
$ws.Range("I7").Value = 1.03979738806215
$ws.Range("J7").Value = 1.048894208843835
$ws.Range("K7").Value = 1.054014393507947
$ws.Range("L7").Value = 1.060444062377748
$ws.Range("M7").Value = 1.067084173896246
$ws.Range("N7").Value = 1.020117691560376
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042853224585336
$ws.Range("D8").Value = 1.050368639143124
$ws.Range("E8").Value = 1.056658180185456
$ws.Range("F8").Value = 1.063298000264248
$ws.Range("I8").Value = 1.03948679188159
$ws.Range("J8").Value = 1.047809194896475
$ws.Range("K8").Value = 1.053060599787044
$ws.Range("L8").Value = 1.059333191942108
$ws.Range("M8").Value = 1.065955354378079
$ws.Range("N8").Value = 1.019748831973483
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040119088004062
$ws.Range("D9").Value = 1.048250504990133
$ws.Range("E9").Value = 1.054272156315677
$ws.Range("F9").Value = 1.060880545171168
$ws.Range("I9").Value = 1.038928746994578
$ws.Range("J9").Value = 1.045892720383075
$ws.Range("K9").Value = 1.051373631175817
$ws.Range("L9").Value = 1.05737625785066
$ws.Range("M9").Value = 1.063964039839361
$ws.Range("N9").Value = 1.019096042727262
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038295873517269
$ws.Range("D10").Value = 1.046837438252679
$ws.Range("E10").Value = 1.052684981709342
$ws.Range("F10").Value = 1.05927054579334
$ws.Range("I10").Value = 1.038549497520085
$ws.Range("J10").Value = 1.044611983346377
$ws.Range("K10").Value = 1.050244757614871
$ws.Range("L10").Value = 1.056072007184128
$ws.Range("M10").Value = 1.062635011075451
$ws.Range("N10").Value = 1.018658954506477
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.0375062740172
$ws.Range("D11").Value = 1.04622532949823
$ws.Range("E11").Value = 1.051998548208311
$ws.Range("F11").Value = 1.058573783003971
$ws.Range("I11").Value = 1.038383565018666
$ws.Range("J11").Value = 1.044056669432194
$ws.Range("K11").Value = 1.049754935425051
$ws.Range("L11").Value = 1.055507338422639
$ws.Range("M11").Value = 1.062059170755968
$ws.Range("N11").Value = 1.018469239450405
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037212959481741
$ws.Range("D12").Value = 1.045997928356798
$ws.Range("E12").Value = 1.051743699749998
$ws.Range("F12").Value = 1.058315030581159
$ws.Range("I12").Value = 1.038321672499894
$ws.Range("J12").Value = 1.043850288173203
$ws.Range("K12").Value = 1.04957284103352
$ws.Range("L12").Value = 1.055297607104644
$ws.Range("M12").Value = 1.061845223189939
$ws.Range("N12").Value = 1.018398702608704
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037275877474593
$ws.Range("D13").Value = 1.046046708340723
$ws.Range("E13").Value = 1.051798360024085
$ws.Range("F13").Value = 1.05837053131333
$ws.Range("I13").Value = 1.038334960323555
$ws.Range("J13").Value = 1.043894562778042
$ws.Range("K13").Value = 1.049611907813199
$ws.Range("L13").Value = 1.055342594654674
$ws.Range("M13").Value = 1.061891118157329
$ws.Range("N13").Value = 1.018413836090643
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037482029007662
$ws.Range("D14").Value = 1.046206533199898
$ws.Range("E14").Value = 1.051977479837044
$ws.Range("F14").Value = 1.05855239330193
$ws.Range("I14").Value = 1.038378454223622
$ws.Range("J14").Value = 1.044039612201748
$ws.Range("K14").Value = 1.049739886562195
$ws.Range("L14").Value = 1.055490001727307
$ws.Range("M14").Value = 1.062041486904477
$ws.Range("N14").Value = 1.018463410245771
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037609042786945
$ws.Range("D15").Value = 1.046305001695696
$ws.Range("E15").Value = 1.05208785782466
$ws.Range("F15").Value = 1.058664451922947
$ws.Range("I15").Value = 1.038405218085684
$ws.Range("J15").Value = 1.044128966947979
$ws.Range("K15").Value = 1.049818718254369
$ws.Range("L15").Value = 1.055580825645424
$ws.Range("M15").Value = 1.062134126773673
$ws.Range("N15").Value = 1.018493945473643
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038348273629187
$ws.Range("D16").Value = 1.046878056771183
$ws.Range("E16").Value = 1.05273055537269
$ws.Range("F16").Value = 1.059316795563631
$ws.Range("I16").Value = 1.038560473768866
$ws.Range("J16").Value = 1.044648821907892
$ws.Range("K16").Value = 1.050277244125307
$ws.Range("L16").Value = 1.056109484060306
$ws.Range("M16").Value = 1.062673220035267
$ws.Range("N16").Value = 1.01867153571195
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038811936067059
$ws.Range("D17").Value = 1.047237454128038
$ws.Range("E17").Value = 1.053133922901225
$ws.Range("F17").Value = 1.059726094281147
$ws.Range("I17").Value = 1.038657402307915
$ws.Range("J17").Value = 1.044974712800047
$ws.Range("K17").Value = 1.050564593861866
$ws.Range("L17").Value = 1.056441118844744
$ws.Range("M17").Value = 1.063011281748905
$ws.Range("N17").Value = 1.018782811984376
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039082369557092
$ws.Range("D18").Value = 1.047447061248985
$ws.Range("E18").Value = 1.053369279834717
$ws.Range("F18").Value = 1.059964867756506
$ws.Range("I18").Value = 1.038713773624149
$ws.Range("J18").Value = 1.045164727489224
$ws.Range("K18").Value = 1.050732102445219
$ws.Range("L18").Value = 1.05663456357811
$ws.Range("M18").Value = 1.063208432519109
$ws.Range("N18").Value = 1.018847673857602
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039174578242696
$ws.Range("D19").Value = 1.047518527924719
$ws.Range("E19").Value = 1.053449543954639
$ws.Range("F19").Value = 1.060046289537962
$ws.Range("I19").Value = 1.038732966725588
$ws.Range("J19").Value = 1.045229505407952
$ws.Range("K19").Value = 1.050789201994472
$ws.Range("L19").Value = 1.056700524560136
$ws.Range("M19").Value = 1.06327564992904
$ws.Range("N19").Value = 1.018869782674074
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038762190829015
$ws.Range("D20").Value = 1.047198896603371
$ws.Range("E20").Value = 1.0530906371635
$ws.Range("F20").Value = 1.059682176634183
$ws.Range("I20").Value = 1.038647019908348
$ws.Range("J20").Value = 1.044939755239863
$ws.Range("K20").Value = 1.050533774059153
$ws.Range("L20").Value = 1.056405536768647
$ws.Range("M20").Value = 1.062975014540064
$ws.Range("N20").Value = 1.018770877614074
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037421323144179
$ws.Range("D21").Value = 1.046159469778682
$ws.Range("E21").Value = 1.051924730108985
$ws.Range("F21").Value = 1.058498837928328
$ws.Range("I21").Value = 1.038365653471796
$ws.Range("J21").Value = 1.043996901886173
$ws.Range("K21").Value = 1.049702204227559
$ws.Range("L21").Value = 1.055446593694867
$ws.Range("M21").Value = 1.061997208573393
$ws.Range("N21").Value = 1.018448813777336
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036578137539882
$ws.Range("D22").Value = 1.045505728363819
$ws.Range("E22").Value = 1.051192392531882
$ws.Range("F22").Value = 1.057755151443189
$ws.Range("I22").Value = 1.038187255605519
$ws.Range("J22").Value = 1.043403438684185
$ws.Range("K22").Value = 1.049178480000237
$ws.Range("L22").Value = 1.054843736268998
$ws.Range("M22").Value = 1.061382105696739
$ws.Range("N22").Value = 1.018245924659812
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037025138983444
$ws.Range("D23").Value = 1.045852309338312
$ws.Range("E23").Value = 1.051580550890916
$ws.Range("F23").Value = 1.058149363006284
$ws.Range("I23").Value = 1.038281969126272
$ws.Range("J23").Value = 1.043718107010217
$ws.Range("K23").Value = 1.049456200015599
$ws.Range("L23").Value = 1.055163316029811
$ws.Range("M23").Value = 1.06170821360563
$ws.Range("N23").Value = 1.018353517499109
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.038784668600396
$ws.Range("D24").Value = 1.047216319160779
$ws.Range("E24").Value = 1.053110195879412
$ws.Range("F24").Value = 1.059702021015556
$ws.Range("I24").Value = 1.038651711778958
$ws.Range("J24").Value = 1.044955551278616
$ws.Range("K24").Value = 1.050547700502891
$ws.Range("L24").Value = 1.056421614753025
$ws.Range("M24").Value = 1.062991402238277
$ws.Range("N24").Value = 1.018776270377229
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040826002303825
$ws.Range("D25").Value = 1.048798264861683
$ws.Range("E25").Value = 1.054888382218776
$ws.Range("F25").Value = 1.06150522578271
$ws.Range("I25").Value = 1.03907428731193
$ws.Range("J25").Value = 1.046388716649252
$ws.Range("K25").Value = 1.051810496891991
$ws.Range("L25").Value = 1.057882106232485
$ws.Range("M25").Value = 1.064479103283156
$ws.Range("N25").Value = 1.019265138092452
